$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D6","D10","D11","D12","D16","D17","D19","D20","D21","D22","D23","D24","D25","D27","D30","D33","D35","D36","D37","D38","D39","D42","D43","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '66.442.60'
$ws.Range('E2').Value = '  +0.68%  '
$ws.Range('D3').Value = '3.297.72'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '587.91'
$ws.Range('E5').Value = '  +2.42%  '
$ws.Range('D6').Value = '180.74'
$ws.Range('E6').Value = '  +1.15%  '
$ws.Range('E7').Value = '  +1.03%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '3.294.81'
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('D10').Value = '0.126'
$ws.Range('E10').Value = '  -0.36%  '
$ws.Range('D11').Value = '6.86'
$ws.Range('E11').Value = '  +2.81%  '
$ws.Range('D12').Value = '0.402'
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('D13').Value = '3.871.99'
$ws.Range('E13').Value = '  -0.05%  '
$ws.Range('E14').Value = '  -2.24%  '
$ws.Range('D15').Value = '66.439.15'
$ws.Range('E15').Value = '  +0.52%  '
$ws.Range('D16').Value = '26.60'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.0000163'
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.284.94'
$ws.Range('E18').Value = '  -0.98%  '
$ws.Range('D19').Value = '426.80'
$ws.Range('E19').Value = '  -2.13%  '
$ws.Range('D20').Value = '5.46'
$ws.Range('E20').Value = '  -2.68%  '
$ws.Range('D21').Value = '13.00'
$ws.Range('E21').Value = '  -2.85%  '
$ws.Range('D22').Value = '7.28'
$ws.Range('E22').Value = '  -2.49%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.28%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '71.49'
$ws.Range('E24').Value = '  -1.87%  '
$ws.Range('D25').Value = '5.73'
$ws.Range('E25').Value = '  +1.16%  '
$ws.Range('D26').Value = '3.448.15'
$ws.Range('E26').Value = '  -0.28%  '
$ws.Range('D27').Value = '0.514'
$ws.Range('E27').Value = '  -0.52%  '
$ws.Range('E28').Value = '  +5.32%  '
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('D30').Value = '9.15'
$ws.Range('E30').Value = '  +1.58%  '
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('E32').Value = '  -1.58%  '
$ws.Range('D33').Value = '22.28'
$ws.Range('E33').Value = '  -1.16%  '
$ws.Range('D35').Value = '5.14'
$ws.Range('E35').Value = '  -0.44%  '
$ws.Range('D36').Value = '6.54'
$ws.Range('E36').Value = '  -2.32%  '
$ws.Range('D37').Value = '1.18'
$ws.Range('E37').Value = '  -1.49%  '
$ws.Range('D38').Value = '159.59'
$ws.Range('E38').Value = '  +1.03%  '
$ws.Range('D39').Value = '1.43'
$ws.Range('E39').Value = '  -2.48%  '
$ws.Range('D40').Value = '2.865.60'
$ws.Range('E40').Value = '  +2.80%  '
$ws.Range('E41').Value = '  -0.64%  '
$ws.Range('D42').Value = '26.23'
$ws.Range('E42').Value = '  -3.58%  '
$ws.Range('D43').Value = '4.31'
$ws.Range('E43').Value = '  -1.34%  '
$ws.Range('E44').Value = '  -4.34%  '
$ws.Range('D45').Value = '39.76'
$ws.Range('E45').Value = '  -1.99%  '
$ws.Range('D46').Value = '0.0657'
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('D47').Value = '5.89'
$ws.Range('E47').Value = '  -3.96%  '
$ws.Range('D48').Value = '2.30'
$ws.Range('E48').Value = '  -0.10%  '
$ws.Range('D49').Value = '312.02'
$ws.Range('E49').Value = '  -3.57%  '
$ws.Range('D50').Value = '22.71'
$ws.Range('E50').Value = '  -3.81%  '
$ws.Range('D51').Value = '0.0270'
$ws.Range('E51').Value = '  -0.58%  '
